$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Exportar Planilha")
$ws2 = $wb.Worksheets.Item("SQL")

# --- Update existing rows 28-43 in "Exportar Planilha" (accumulated totals shift by one month) ---
$ws1.Range("H28").Value2 = 17.0
$ws1.Range("H29").Value2 = 15.0
$ws1.Range("H30").Value2 = 14.0
$ws1.Range("E31").Value2 = 5353535.0
$ws1.Range("F31").Value2 = 117014383.64
$ws1.Range("H31").Value2 = 14.0
$ws1.Range("H32").Value2 = 12.0
$ws1.Range("H33").Value2 = 12.0
$ws1.Range("E34").Value2 = 6092265.0
$ws1.Range("F34").Value2 = 125086128.42
$ws1.Range("H34").Value2 = 11.0
$ws1.Range("E35").Value2 = 6204170.0
$ws1.Range("F35").Value2 = 123715934.42
$ws1.Range("H35").Value2 = 11.0
$ws1.Range("E36").Value2 = 5718034.0
$ws1.Range("F36").Value2 = 123955303.99
$ws1.Range("H36").Value2 = 9.0
$ws1.Range("E37").Value2 = 5690223.0
$ws1.Range("F37").Value2 = 114018670.4
$ws1.Range("H37").Value2 = 10.0
$ws1.Range("E38").Value2 = 6398850.0
$ws1.Range("F38").Value2 = 130633886.88
$ws1.Range("H38").Value2 = 7.0
$ws1.Range("E39").Value2 = 5379442.0
$ws1.Range("F39").Value2 = 107995271.46
$ws1.Range("H39").Value2 = 6.0
$ws1.Range("E40").Value2 = 6817527.0
$ws1.Range("F40").Value2 = 141200864.72
$ws1.Range("H40").Value2 = 5.0
$ws1.Range("E41").Value2 = 6716608.0
$ws1.Range("F41").Value2 = 125895727.56
$ws1.Range("H41").Value2 = 4.0
$ws1.Range("E42").Value2 = 7123312.0
$ws1.Range("F42").Value2 = 142552015.45
$ws1.Range("H42").Value2 = 3.0
$ws1.Range("E43").Value2 = 6199778.0
$ws1.Range("F43").Value2 = 117811130.85
$ws1.Range("H43").Value2 = 2.0

# --- Add new row 44 (201707) — copy formatting from row 43 first, then fill values ---
$ws1.Range("A43:H43").Copy()
$ws1.Range("A44:H44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A44").Value2 = "'201707"
$ws1.Range("B44").Value2 = 6967247.0
$ws1.Range("C44").Value2 = 126904132.21
$ws1.Range("D44").Value2 = 2582971.22
$ws1.Range("E44").Value2 = 2828194.0
$ws1.Range("F44").Value2 = 69248014.87
$ws1.Range("G44").Value2 = 2582971.22
$ws1.Range("H44").Value2 = 2.0

# Re-paste formats from row 43 so the leading apostrophe on A44 does not leave a stray quote-prefix style
$ws1.Range("A43:H43").Copy()
$ws1.Range("A44:H44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update SQL text on "SQL" sheet: extend date range through 201707 ---
$sql = @"
select  t.anomes,
        sum(nvl(t.QTA_FT,0)) QTA_FT,
        sum(nvl(t.VA_FT,0)) VA_FT,
        sum(nvl(t.VR_FT,0)) VR_FT,        
        sum(nvl(t.QTA_AT,0)) QTA_AT,
        sum(nvl(t.VA_AT,0)) VA_AT,
        sum(nvl(t.VR_FT,0)) VR_AT,        
        sum(t.C_FAT) C_AFT        
from
(select substr(fc.ID_TEMPO_MES_ANO_REF,1,6) anomes,
        substr(fc.ID_TEMPO_MES_ANO_REF,1,6) anomes_f,
        0 C_FAT,        
        sum(nvl(fc.QTD_ITEM,0)) - sum(nvl(fc.QTD_GLOSADO,0)) QTA_FT,
        sum(nvl(fc.VAL_APROVADO_ITEM,0)) VA_FT,
        sum(nvl(fc.VALOR_PAGO_REVISAO,0)) VR_FT,        
        count(1) N_FT,
        sum(0) QTA_AT,
        sum(0) VA_AT,
        sum(0) VR_AT,        
        sum(0) N_AT
from    TS.FAT_ITEM_CONTA fc
group by substr(fc.ID_TEMPO_MES_ANO_REF,1,6),
         0
union   
select  to_char(DATA_ATENDIMENTO, 'RRRRMM') anomes,
        substr(fc.ID_TEMPO_MES_ANO_REF,1,6) anomes_f,
        1 C_FAT,        
        sum(0) QTA_FT,
        sum(0) VA_FT,
        sum(0) VR_FT,        
        sum(0) N_FT,
        sum(nvl(fc.QTD_ITEM,0)) - sum(nvl(fc.QTD_GLOSADO,0)) QTA_AT,
        sum(nvl(fc.VAL_APROVADO_ITEM,0)) VA_AT,
        sum(nvl(fc.VALOR_PAGO_REVISAO,0)) VR_AT,        
        count(1) N_AT
from    TS.FAT_ITEM_CONTA fc
group by to_char(DATA_ATENDIMENTO, 'RRRRMM'),
         substr(fc.ID_TEMPO_MES_ANO_REF,1,6),
         1) t
where    t.anomes between '201401' and '201707'         
group by t.anomes
having  sum(nvl(t.QTA_FT,0)) + sum(nvl(t.VA_FT,0)) + sum(nvl(t.QTA_AT,0)) + sum(nvl(t.VA_AT,0)) <> 0
and     substr(t.anomes,1,3) = '201'
order by 1
"@
$ws2.Range("A2").Value2 = $sql
